$wb = $excel.ActiveWorkbook

# --- dailyQuests sheet: rework the "key" column into a 0-based "index" column ---
$wsQuests = $wb.Worksheets.Item("dailyQuests")

# Header: STR_key -> INT_index
$wsQuests.Range("A1").Value = "INT_index"

# Column A (previously duplicated the quest_N key, now a literal 0-based index)
$wsQuests.Range("A2").Value = 0
$wsQuests.Range("A3").Value = 1
$wsQuests.Range("A4").Value = 2
$wsQuests.Range("A5").Value = 3
$wsQuests.Range("A6").Value = 4
$wsQuests.Range("A7").Value = 5
$wsQuests.Range("A8").Value = 6
$wsQuests.Range("A9").Value = 7
$wsQuests.Range("A10").Value = 8
$wsQuests.Range("A11").Value = 9

# --- dailyQuestStar sheet: selection grows to cover the whole column ---
$wsStar = $wb.Worksheets.Item("dailyQuestStar")
$wsStar.Range("B2:B6").Select()

# --- dailyQuestStyle sheet: no longer the active tab, selection resets to A2 ---
$wsStyle = $wb.Worksheets.Item("dailyQuestStyle")
$wsStyle.Range("A2").Select()

# This sheet becomes the active / selected tab, with A2 selected (must be
# activated last so it ends up as the workbook's active sheet)
$wsQuests.Activate()
$wsQuests.Range("A2").Select()
